$d = $word.ActiveDocument

# 1. Add a "Working: " lead-in run before "Make navbar format nicer"
$findRange = $d.Content.Duplicate
$findRange.Find.Execute("Make navbar format nicer", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$findRange.Collapse(1)
$findRange.InsertBefore("Working: ")

# 2. Remove the "Put labels and heading on charts" paragraph entirely
$delRange = $d.Content.Duplicate
$delRange.Find.Execute("Put labels and heading on charts", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$delRange.Expand(4)
$delRange.Delete()

# 3. Remove the "Load detail of watchlist including company summary statistics?" paragraph entirely
$delRange2 = $d.Content.Duplicate
$delRange2.Find.Execute("Load detail of watchlist including company summary statistics?", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$delRange2.Expand(4)
$delRange2.Delete()
